$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '63.394.96'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '3.250.30'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '594.77'
$ws.Range("D6").Value = '140.92'
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.247.53'
$ws.Range("E8").Value = '  +3.73%  '
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("D11").Value = '5.40'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").Value = '0.466'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").Value = '34.28'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '3.783.37'
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '3.241.13'
$ws.Range("E17").Value = '  +3.51%  '
$ws.Range("D18").Value = '63.422.25'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '474.26'
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("D22").Value = '0.732'
$ws.Range("E22").Value = '  +3.98%  '
$ws.Range("D23").Value = '7.97'
$ws.Range("E23").Value = '  +3.80%  '
$ws.Range("D24").Value = '83.51'
$ws.Range("E24").Value = '  -4.67%  '
$ws.Range("D25").Value = '13.19'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '7.24'
$ws.Range("E28").Value = '  +4.24%  '
$ws.Range("D29").Value = '8.07'
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("E30").Value = '  +4.06%  '
$ws.Range("D31").Value = '27.69'
$ws.Range("E31").Value = '  +1.76%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("E34").Value = '  -2.25%  '
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("D37").Value = '52.67'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = '0.0₃0713'
$ws.Range("E38").Value = '  -2.56%  '
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '420.85'
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("D41").Value = '2.989.57'
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").Value = '8.36'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("D43").Value = '2.75'
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("D44").Value = '0.109'
$ws.Range("E44").Value = '  -8.03%  '
$ws.Range("E45").Value = '  +2.74%  '
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D48").Value = '25.85'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").Value = '2.32'
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '122.85'
$ws.Range("E51").Value = '  +1.94%  '
